# Rename the original sheet from "Sheet1" to "Data"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Data"

# Zoom the Data sheet view to 120%
$ws1.Activate()
$excel.ActiveWindow.Zoom = 120

# Add a new worksheet right after "Data" and name it "Codebook"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Codebook"

# --- Title row ---
$ws2.Range("A1").Value = "Codebook - explanations of each variable in dataset"
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Size = 14

# --- Header row (row 3) ---
$ws2.Range("B3").Value = "Variable Name"
$ws2.Range("C3").Value = "Variable Label"
$ws2.Range("D3").Value = "Values Defined (if applicable)"
$ws2.Range("B3:D3").Font.Bold = $true
$ws2.Range("B3:D3").Interior.Color = 65535
$ws2.Range("B3:D3").Borders.Item(9).LineStyle = 1

# --- Variable rows ---
$ws2.Range("B4").Value = "SubjectID"
$ws2.Range("C4").Value = "Subject ID"

$ws2.Range("B5").Value = "Age"
$ws2.Range("C5").Value = "Age in Years"

$ws2.Range("B6").Value = "WeightPRE"
$ws2.Range("C6").Value = "Weight in Pounds - Before Program"

$ws2.Range("B7").Value = "WeightPOST"
$ws2.Range("C7").Value = "Weight in Pounds - After Program"

$ws2.Range("B8").Value = "Height"
$ws2.Range("C8").Value = "Height in Decimal Feet"

$ws2.Range("B9").Value = "SES"
$ws2.Range("C9").Value = "Pseudo Socio-Economic-Status"
$ws2.Range("D9").Value = "1=low income; 2=average income; 3=high income"

$ws2.Range("B10").Value = "GenderSTR"
$ws2.Range("C10").Value = "Gender as a Character/Text"

$ws2.Range("B11").Value = "GenderCoded"
$ws2.Range("C11").Value = "Gender Recoded"
$ws2.Range("D11").Value = "1=Male; 2=Female"

$qtext = "1=none of the time; 2=a little of the time; 3=some of the time; 4=a lot of the time; 5=all of the time"

$ws2.Range("B12").Value = "q1"
$ws2.Range("C12").Value = "Hypothetical Question 1"
$ws2.Range("D12").Value = $qtext

$ws2.Range("B13").Value = "q2"
$ws2.Range("C13").Value = "Hypothetical Question 2"
$ws2.Range("D13").Value = $qtext

$ws2.Range("B14").Value = "q3"
$ws2.Range("C14").Value = "Hypothetical Question 3"
$ws2.Range("D14").Value = $qtext

$ws2.Range("B15").Value = "q4"
$ws2.Range("C15").Value = "Hypothetical Question 4"
$ws2.Range("D15").Value = $qtext

$ws2.Range("B16").Value = "q5"
$ws2.Range("C16").Value = "Hypothetical Question 5"
$ws2.Range("D16").Value = $qtext

$ws2.Range("B17").Value = "q6"
$ws2.Range("C17").Value = "Hypothetical Question 6"
$ws2.Range("D17").Value = $qtext

$ws2.Range("D12:D17").WrapText = $true

# Column widths roughly matching the authored layout
$ws2.Columns.Item(2).ColumnWidth = 14.6328125
$ws2.Columns.Item(3).ColumnWidth = 30.36328125
$ws2.Columns.Item(4).ColumnWidth = 45.08984375

# Zoom the Codebook sheet view to 120% as well
$ws2.Activate()
$excel.ActiveWindow.Zoom = 120

# Re-activate the Data sheet (matches tabSelected="1" on Data in the source)
$ws1.Activate()

Write-Output "done"
